# Add a new row (69) to the roll-check rules sheet:
#   A69 = item name
#   D69 = unit rule (local)
# This pushes the sheet's used range from A1:H68 to A1:H69 and Excel
# recalculates row spans / dimension automatically when the cells are
# written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "W3 SCC #7 0230驅動ROLLER 舊品拆裝"
$ws.Range("D69").Value = "1SET=2PC"

# Move the active selection to the newly added row, mirroring the
# author moving to/filling in the new last row before saving.
$ws.Range("E69").Select() | Out-Null
